$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D price values are stored as TEXT (not numbers) in this sheet. ---
# Assigning a numeric-looking string straight to .Value would make Excel
# auto-convert it to a real number (and lose formatting like trailing zeros,
# e.g. "0.1430" -> 0.143), so each of these cells is temporarily switched to
# a Text number format before the value is written, then the cell style is
# restored to "Normal" (its original, unstyled state) afterward.
$priceUpdates = @{
    "D2" = "246.98"
    "D4" = "5.446"
    "D5" = "0.05664"
    "D7" = "0.8011"
    "D8" = "1.031"
    "D9" = "0.01157"
    "D10" = "0.1430"
    "D11" = "0.07238"
    "D12" = "0.03159"
    "D13" = "0.02952"
    "D14" = "0.09281"
    "D15" = "0.001628"
    "D16" = "3.226"
    "D17" = "0.04730"
    "D18" = "0.006478"
    "D19" = "0.005015"
    "D21" = "0.0001502"
    "D22" = "0.0003203"
    "D23" = "3.854"
    "D25" = "2.089"
    "D26" = "0.3276"
    "D27" = "0.1298"
    "D40" = "0.04080"
    "D41" = "0.006949"
    "D42" = "0.1039"
    "D43" = "0.003204"
    "D44" = "0.009036"
    "D45" = "0.00005850"
    "D47" = "0.7860"
    "D48" = "0.009838"
    "D49" = "0.00002102"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# --- Plain text fields (coin name, coinranking link, volume label). ---
# These never look like numbers, so Excel keeps them as text automatically.
$textUpdates = @{
    "B9" = "One"
    "C9" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "E9" = "8OneONEBestin24h"
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "E10" = "9WazirXWRX"
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "E11" = "10MandalaExchangeTokenMDX"
    "B12" = "LiechtensteinCryptoassetsExchange"
    "C12" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "E12" = "11LiechtensteinCryptoassetsExchangeLCX"
    "B13" = "BitrueCoin"
    "C13" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "E13" = "12BitrueCoinBTR"
    "B14" = "BitMartToken"
    "C14" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "E14" = "13BitMartTokenBMX"
    "B15" = "BitForexToken"
    "C15" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "E15" = "14BitForexTokenBF"
    "B16" = "MCDex"
    "C16" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "E16" = "15MCDexMCB"
    "B17" = "CoinExToken"
    "C17" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "E17" = "16CoinExTokenCET"
    "E19" = "18HotbitTokenHTB"
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}
